$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Insert 3 new blank rows just above the "Total" row (old row 17),
#    one at a time so each push-down keeps the style of the rows it
#    displaces (Total row / footer row keep their own formatting).
# ------------------------------------------------------------------
$ws.Rows.Item(17).Insert()
$ws.Rows.Item(18).Insert()
$ws.Rows.Item(19).Insert()

# Give the 3 new rows (17:19) the same cell formatting as the last
# product row (row 16) by copying formats only.
$ws.Range("A16:Q16").Copy()
$ws.Range("A17:Q19").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Match the row heights used in the final layout.
$ws.Rows.Item(17).RowHeight = 25.5
$ws.Rows.Item(18).RowHeight = 24.75
$ws.Rows.Item(19).RowHeight = 25.5
$ws.Rows.Item(20).RowHeight = 24.75

# ------------------------------------------------------------------
# 2. (Re)write every product row, 7 through 19, with the final,
#    alphabetically-sorted list of items (this both updates existing
#    rows and fills in the 3 newly inserted ones).
# ------------------------------------------------------------------
$items = @(
    @("BLOCK AND WHITE  CREAM",               "0:0", "0", "139.00", "139.0000", "1:0"),
    @("CETAL COLD & FLU 20 CAPLETS",           "1:0", "1", "36.00",  "18.0000",  "0:1"),
    @("CORDO PLUS SPRAY 60 ML",                "0:0", "1", "55.00",  "55.0000",  "1:0"),
    @("DECLOPHEN 75MG/3ML 3 AMPOULES",         "5:2", "1", "36.00",  "11.8800",  "0:1"),
    @("EPICEPHIN 2GM I.V. VIAL",               "3:0", "1", "77.00",  "154.0000", "2:0"),
    @("FLAGYL 500MG 20 TAB.",                  "2:0", "1", "34.00",  "17.0000",  "0:1"),
    @("FLOXAMO 500/500 MG 16 F.C.TABS",        "0:1", "1", "110.00", "55.0000",  "0:1"),
    @("FLUMOX 500MG 16 CAPS",                  "0:0", "1", "71.00",  "35.5000",  "0:1"),
    @("FUSI 2% OINT. 15 GM",                   "2:0", "1", "35.00",  "35.0000",  "1:0"),
    @("VASTAFLAM 50MG 20 SUGAR COATED TAB.",   "1:1", "1", "36.00",  "18.0000",  "0:1"),
    @([char]0x062C + [char]0x0647 + [char]0x0627 + [char]0x0632 + " " + [char]0x0631 + [char]0x064A + [char]0x062F, "6:0", "0", "140.00", "140.0000", "1:0"),
    @([char]0x0633 + [char]0x0631 + [char]0x0646 + [char]0x062C + [char]0x0627 + [char]0x062A + " 10 " + [char]0x0633 + [char]0x0645, "0:0", "0", "4.00", "8.0000", "2:0"),
    @([char]0x0633 + [char]0x0631 + [char]0x0646 + [char]0x062C + [char]0x0627 + [char]0x062A + " 3 " + [char]0x0633 + [char]0x0645, "0:0", "0", "2.00", "4.0000", "2:0")
)

$row = 7
foreach ($item in $items) {
    $ws.Range("A" + $row).Value = ($row - 6)
    $ws.Range("C" + $row).Value = $item[0]
    $ws.Range("H" + $row).Value = $item[1]
    $ws.Range("L" + $row).Value = $item[2]
    $ws.Range("N" + $row).Value = $item[3]
    $ws.Range("P" + $row).Value = $item[4]
    $ws.Range("Q" + $row).Value = $item[5]
    $row = $row + 1
}

# ------------------------------------------------------------------
# 3. Update the total (sum of sell prices) and the generated-on
#    timestamp in the footer.
# ------------------------------------------------------------------
$ws.Range("P20").Value = 690.38
$ws.Range("K21").Value = "Wednesday, 4 June, 2025 10:19 AM"

Write-Host "Edit applied"
